# Update financial figures on the CMTL worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CMTL")

# Row 21 - Earnings Before Interest And Taxes
$ws.Range("D21").Value = 69600
$ws.Range("E21").Value = 74300
$ws.Range("F21").Value = 22800
$ws.Range("G21").Value = 47200
$ws.Range("H21").Value = 57800
$ws.Range("I21").Value = 49800
$ws.Range("J21").Value = 69700

# Row 24 - Income Tax Expense
$ws.Range("D24").Value = 6200

# Row 26 - Income After Tax
$ws.Range("D26").Value = 18500

# Row 27 - Net Income From Continuing Ops
$ws.Range("D27").Value = 18500

# Row 29 - Discontinued Operations
$ws.Range("D29").Value = 11300

# Row 89 - Total Cash Flow From Operating Activities
$ws.Range("E89").Value = 66900
$ws.Range("F89").Value = 15100

# Row 100 - Total Cash Flows From Financing Activities
$ws.Range("E100").Value = -83700
$ws.Range("F100").Value = 187000
